$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 13.00385802469133
$ws.Range("N2").Value = 1.587252942212933
$ws.Range("O2").Value = 1.701551834435819

$ws.Range("I3").Value = 13.00385802469133

$ws.Range("I4").Value = 18.89814814814816

$ws.Range("I5").Value = -1.819444444444444
$ws.Range("N5").Value = 1.421349525022035
$ws.Range("O5").Value = 1.510209672321146

$ws.Range("I6").Value = -1.819444444444444
